$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("heartbeat")
Write-Host "Orientation before:" $ws2.PageSetup.Orientation()
$ws2.PageSetup.Orientation = 1
Write-Host "Orientation after:" $ws2.PageSetup.Orientation()
